# Insert a new data row at row 46 (pushing the existing rows 46-94 down to 47-95)
# and populate it with the new Alcachofa record for Macroferia Regional de Talca.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(46).Insert()

$ws.Cells(46,1).Value  = 5
$ws.Cells(46,2).Value  = "Macroferia Regional de Talca"
$ws.Cells(46,3).Value  = "Maule"
$ws.Cells(46,4).Value  = 44789
$ws.Cells(46,5).Value  = 7
$ws.Cells(46,6).Value  = 100112013
$ws.Cells(46,7).Value  = "Alcachofa"
$ws.Cells(46,8).Value  = "Madrigal"
$ws.Cells(46,9).Value  = "Primera"
$ws.Cells(46,10).Value = 400
$ws.Cells(46,11).Value = 14000
$ws.Cells(46,12).Value = 14000
$ws.Cells(46,13).Value = 14000
$ws.Cells(46,14).Value = "`$/caja 40 unidades"
$ws.Cells(46,15).Value = "Provincia del Elquí"
$ws.Cells(46,16).Value = 350
$ws.Cells(46,17).Value = 40
$ws.Cells(46,18).Value = "Hortaliza"
